# add NGSCheckMate to S4
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("workflow_repository_tools")

# Insert a new row above row 2 (shifts existing rows 2..20 down to 3..21)
$ws.Rows.Item(2).Insert()

# Populate the newly inserted row 2 with the NGSCheckMate entry.
# The version "1.3" must land as a shared text string (not a number),
# matching how the existing version-number column is stored. Setting
# .Value = "1.3" directly would get auto-coerced to a numeric cell, and
# forcing text via NumberFormat="@" (or a leading quote) would mint a new,
# unused cell-style entry that isn't part of this change. Instead, copy the
# value from a cell elsewhere in the workbook that already holds the exact
# literal text "1.3" (found via Find), which carries its text-ness over
# without touching styles.
$ws.Cells.Item(2, 1).Value = "NGSCheckMate"
$src = $wb.Worksheets.Item("r_packages")
$textOneThree = $src.Columns.Item(2).Find("1.3", [Type]::Missing, -4123, 1)
$textOneThree.Copy($ws.Cells.Item(2, 2))
$ws.Cells.Item(2, 3).Value = "https://github.com/d3b-center/OpenPBTA-workflows/blob/master/cwl/bcf_call.cwl"
